$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) AMSIN sheet - normalize formatting on the latest registration row
#    (row 19) so it matches the plain/default style used by the rows
#    above it. The B column keeps its existing date/time number
#    format, everything else goes back to the default "Normal" style.
# ------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$wsAmsin.Range("A19").Style = "Normal"
$wsAmsin.Range("C19:G19").Style = "Normal"

# Run time got re-saved with full precision.
$wsAmsin.Range("B19").Value = 44810.94112326389

# ------------------------------------------------------------------
# 2) AMS sheet - append the newest registration-script run that came
#    in from the live environment.
# ------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Pull the Run Time number format (yyyy-mm-dd hh:mm:ss) from an
# existing timestamp cell so the new row matches the rest of the data.
$wsAmsin.Range("B18").Copy() | Out-Null
$wsAms.Range("B10").PasteSpecial(-4122) | Out-Null

$wsAms.Range("A10").Value = "'2022-09-08"
$wsAms.Range("B10").Value = 44812.52689823115
$wsAms.Range("C10").Value = "educ166"
$wsAms.Range("D10").Value = 60
$wsAms.Range("E10").Value = 60
$wsAms.Range("F10").Value = 0
$wsAms.Range("G10").Value = 1.06
